# Auto-generated edit script applying the cryptos.xlsx price/volume refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.112.14'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.841.77'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.44'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6859'
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3021'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07448'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.12'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07668'
$ws.Range('D12').Value = '1.844.31'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.061'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6833'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '87.52'
$ws.Range('E15').Value = '  -5.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.171'
$ws.Range('E16').Value = '  -6.98%  '
$ws.Range('D17').Value = '29.111.15'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008158'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('D19').Value = '2.081.11'
$ws.Range('E19').Value = '  -0.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '227.81'
$ws.Range('E20').Value = '  -5.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.387'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.27'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1452'
$ws.Range('E26').Value = '  -3.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.763'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.10'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.510'
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.269'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.137'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.194'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05230'
$ws.Range('E33').Value = '  +2.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7609'
$ws.Range('E34').Value = '  -4.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.851'
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.134'
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.689'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = '1.309.50'
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01839'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9327'
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.807'
$ws.Range('E42').Value = '  -3.94%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.05'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9998'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000123'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.983.78'
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5200'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '64.82'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.522'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.772'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('B51').Value = 'XinFinNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07374'
$ws.Range('E51').Value = '  +16.84%  '
